$wb = $excel.ActiveWorkbook

# --- Rename the only sheet: "regData" -> "Sheet1" ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Sheet1"

# --- Move the cell selection/active cell on that sheet from A4 to B11 ---
$ws.Activate()
$ws.Range("B11").Select()

# --- Reposition / resize the workbook window (xWindow/yWindow/windowWidth/windowHeight) ---
$win = $wb.Windows.Item(1)
$win.Left = 4800
$win.Top = 1620
$win.Width = 14400
$win.Height = 8170
